# hosteller logic updated for report
# Replace the absentee row with the second student's details and remove the
# now-duplicate row, shrinking the table from two student rows to one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 to hold the second student's info (22ALR056 / NAVEEN SAKTHI S)
$ws.Range("B5").Value = "22ALR056"
$ws.Range("C5").Value = "NAVEEN SAKTHI S"

# Row 4 (header row) height changes from 20 to 25
$ws.Rows.Item(4).RowHeight = 25

# Delete the now redundant row 6 entirely
$ws.Rows.Item(6).Delete()
